$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 8 ("Display an alert on the article"): merge the trailing runs of
#    the caption textbox into a single run, keeping the first two runs as-is.
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$caption = $slide8.Shapes.Item(4).TextFrame.TextRange
$merged = $caption.Characters(38, 91)
$merged.Text = " datasets associated with both articles in the references AND with the current publication "

# ---------------------------------------------------------------------------
# 2. Append a new slide ("Our code") with a github link + license note.
# ---------------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Our code"
$title.LanguageID = "en-GB"

$url = "https://github.com/ScientificDataLabs/dataset-updates-plugin"
$urlPart1 = "https://"
$urlPart2 = "github.com/ScientificDataLabs/dataset-updates-plugin"

$content = $newSlide.Shapes.Item(2).TextFrame.TextRange
$content.Text = $url
$content.Font.Size = 24
$content.LanguageID = "en-GB"
$content.InsertAfter("`rMIT licensed")

$link1 = $content.Characters(1, $urlPart1.Length)
$link1.ActionSettings(1).Hyperlink.Address = $url

$link2 = $content.Characters($urlPart1.Length + 1, $urlPart2.Length)
$link2.ActionSettings(1).Hyperlink.Address = $url
